$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Offers"
$ws.Range("B7").Value = "offerID(int)"
$ws.Range("B7").Font.Underline = $true
$ws.Range("C7").Value = "farmerID(int) REFERENCES farmers(farmerID)"
$ws.Range("C7").Font.Underline = $true
$ws.Range("C7").HorizontalAlignment = -4108
$ws.Range("C7").WrapText = $true

$ws.Range("A10").Value = "Products"
$ws.Range("B11").Value = "productID(int)"
$ws.Range("B11").Font.Underline = $true
$ws.Range("C11").Value = "productName(varchar(50))"
$ws.Range("C11").WrapText = $true

$ws.Range("A14").Value = "ProductsInOffer"
$ws.Range("C15").Value = "productID(int) REFERENCES Products(productID) "
$ws.Range("B15").Value = "offerID(int) REFERENCES Offers(offerID)"
$ws.Range("B15").Font.Underline = $true
$ws.Range("B15").WrapText = $true
$ws.Range("C15").Font.Underline = $true
$ws.Range("C15").WrapText = $true

$ws.Rows(7).AutoFit()
$ws.Rows(11).AutoFit()
$ws.Rows(15).AutoFit()
